# update script variable Es Adulto OK
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column V = "Edad" (Age), column X = "Indigena", column AG = "Taller - Actividad"
# Row 2 (first data record): age corrected 25 -> 17, "Indigena" flag 4 -> 1,
# and the "Taller - Actividad" cell was a stray text value "1a" that is now a
# proper number 1.
$ws.Range("V2").Value = 17
$ws.Range("X2").Value = 1
$ws.Range("AG2").Value = 1

# Row 3 (second data record): age corrected 1 -> 17
$ws.Range("V3").Value = 17

# Scroll the view and move the active selection to match where the editor
# ended up after making the change.
$excel.ActiveWindow.ScrollColumn = 16
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("V4").Select() | Out-Null
